# Refactoring the OngyAi app
# Updates the Consumption Forecast data: new forecasted consumption values
# (column A) and timestamps shifted forward by 46 days (column B), for
# rows 2 through 97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(5090,5060,5030,5010,5000,5000,5000,5000,5000,5000,5000,5010,5040,5070,5120,5180,5240,5350,5480,5650,5850,6030,6230,6410,6610,6780,6930,7070,7210,7300,7370,7400,7440,7440,7430,7400,7350,7310,7260,7200,7150,7130,7110,7080,7050,7030,7020,7000,7020,7020,7020,7020,7020,7020,7020,7020,7030,7030,7030,7030,7050,7070,7110,7170,7220,7290,7360,7430,7520,7620,7700,7800,7930,8010,8060,8050,8010,7950,7870,7790,7670,7530,7420,7290,7100,6960,6810,6660,6510,6370,6260,6150,6100,6040,5990,5930)

$dayShift = 46

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
    $ws.Cells.Item($row, 2).Value = $ws.Cells.Item($row, 2).Value2 + $dayShift
}
